# Slide 7 ("RStudio cont."), Content Placeholder 2: remove the three
# sub-bullets "Highlighting a word", "Changing the theme", and
# "Using the help tab" that used to follow "How to open a script".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 7 is "Highlighting a word" before any deletions; each
# .Delete() shifts the following paragraphs up by one, so re-fetching
# paragraph 7 three times removes the three target bullets in order
# while leaving the trailing empty paragraph untouched.
for ($i = 0; $i -lt 3; $i++) {
    $target = $tr.Paragraphs(7, 1)
    $target.Delete()
}
